{"js": "// Insert the new intro paragraphs (and the \"3 Creating Objects\" heading)\n// right after the \"Write Up\" title paragraph, before the existing blank\n// paragraph that used to follow it.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The paragraph immediately after the \"Write Up\" title is the empty\n// paragraph that originally followed it; anchor the new content there so\n// it lands between the title and that blank paragraph (which is left\n// untouched), matching the target layout.\nconst anchor = paragraphs.items[1];\n\nconst introText =\n  \"This week, we will be taking a look at how simple it is to create an \" +\n  \"object in Game Maker. You have created your sprite in last\\u2019s weeks \" +\n  \"lesson, but that is only the face of your game elements. In order to \" +\n  \"get anything, you place into a game, to do anything. You are going to \" +\n  \"have to attach that face to an object, and I will explain how to go \" +\n  \"about doing that in this tutorial here.\";\n\nconst calloutText =\n  \"So, if this at all sounds interesting to you, then please join us for \" +\n  \"this week\\u2019s article entitled: \";\n\nconst headingText = \"3 Creating Objects\";\n\nanchor.insertParagraph(introText, \"Before\");\nanchor.insertParagraph(calloutText, \"Before\");\nconst headingParagraph = anchor.insertParagraph(headingText, \"Before\");\nheadingParagraph.style = \"Heading 1\";\n\nawait context.sync();\n", "ps1": "# Insert the new intro paragraphs (and the \"3 Creating Objects\" heading)\n# right after the \"Write Up\" title paragraph, before the existing blank\n# paragraph that used to follow it.\n\n$d = $word.ActiveDocument\n\n$introText = \"This week, we will be taking a look at how simple it is to create an object in Game Maker. You have created your sprite in last\" + [char]0x2019 + \"s weeks lesson, but that is only the face of your game elements. In order to get anything, you place into a game, to do anything. You are going to have to attach that face to an object, and I will explain how to go about doing that in this tutorial here.\"\n$calloutText = \"So, if this at all sounds interesting to you, then please join us for this week\" + [char]0x2019 + \"s article entitled: \"\n$headingText = \"3 Creating Objects\"\n\n# Paragraph 2 is the blank paragraph that originally followed \"Write Up\";\n# it stays in place and we insert the three new paragraphs before it\n# (re-fetching by index each time since paragraph handles go stale after\n# a structural edit).\n$anchorIndex = 2\n\n$d.Paragraphs.Item($anchorIndex).Range.InsertParagraphBefore()\n$d.Paragraphs.Item($anchorIndex).Range.Text = $introText\n\n$d.Paragraphs.Item($anchorIndex + 1).Range.InsertParagraphBefore()\n$d.Paragraphs.Item($anchorIndex + 1).Range.Text = $calloutText\n\n$d.Paragraphs.Item($anchorIndex + 2).Range.InsertParagraphBefore()\n$d.Paragraphs.Item($anchorIndex + 2).Range.Text = $headingText\n$d.Paragraphs.Item($anchorIndex + 2).Style = \"Heading 1\"\n"}
